$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 updates ---
# B2: company_name placeholder changes from "2" to "1" (kept as text, not number)
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "1"
$ws.Range("B2").Style = "Normal"

$ws.Range("D2").Value = 0.177
$ws.Range("E2").Value = -0.208
$ws.Range("G2").Value = 1.774436090225564
$ws.Range("H2").Value = 1.774436090225564
$ws.Range("I2").Value = 0.08721804511278196
$ws.Range("J2").Value = 0.08721804511278196
$ws.Range("K2").Value = 0.079
$ws.Range("L2").Value = 0.118796992481203
$ws.Range("M2").Value = 0.237
$ws.Range("N2").Value = 0.03726415094339623
$ws.Range("O2").Value = 3
$ws.Range("P2").Value = 0.237
$ws.Range("Q2").Value = 0.03726415094339623
$ws.Range("R2").Value = 3
$ws.Range("U2").Value = 0
$ws.Range("V2").Value = 0
$ws.Range("W2").Value = 0.02358208955223881
$ws.Range("X2").Value = 0.0473054834734378
$ws.Range("Y2").Value = -0.02372339392119899
$ws.Range("Z2").Value = 0.1985074626865672
$ws.Range("AA2").Value = 0.0173134328358209
$ws.Range("AB2").Value = 0.04712209044112428
$ws.Range("AC2").Value = -0.02980865760530339
$ws.Range("AD2").Value = 0.07199999999999999
$ws.Range("AF2").Value = 0.07199999999999999
$ws.Range("AG2").Value = 0.07199999999999999
$ws.Range("AH2").Value = 0.01119402985074627
$ws.Range("AI2").Value = 0.02635431918008784
$ws.Range("AJ2").Value = 0.01119402985074627
$ws.Range("AK2").Value = 0.02635431918008784
$ws.Range("AL2").Value = 0.008
$ws.Range("AM2").Value = 0.008
$ws.Range("AN2").Value = 1.241379310344827
$ws.Range("AO2").Value = 7.25
$ws.Range("AP2").Value = 1.241379310344827
$ws.Range("AQ2").Value = 7.25

# --- Row 3 updates ---
$ws.Range("D3").Value = 0.177
$ws.Range("E3").Value = -0.208
$ws.Range("G3").Value = 1.774436090225564
$ws.Range("H3").Value = 1.774436090225564
$ws.Range("I3").Value = 0.08721804511278196
$ws.Range("J3").Value = 0.08721804511278196
$ws.Range("K3").Value = 0.079
$ws.Range("L3").Value = 0.118796992481203
$ws.Range("M3").Value = 0.237
$ws.Range("N3").Value = 0.03726415094339623
$ws.Range("O3").Value = 3
$ws.Range("P3").Value = 0.237
$ws.Range("Q3").Value = 0.03726415094339623
$ws.Range("R3").Value = 3
$ws.Range("U3").Value = 0
$ws.Range("V3").Value = 0
$ws.Range("W3").Value = 0.02358208955223881
$ws.Range("X3").Value = 0.0473054834734378
$ws.Range("Y3").Value = -0.02372339392119899
$ws.Range("Z3").Value = 0.1985074626865672
$ws.Range("AA3").Value = 0.0173134328358209
$ws.Range("AB3").Value = 0.04712209044112428
$ws.Range("AC3").Value = -0.02980865760530339
$ws.Range("AD3").Value = 0.07199999999999999
$ws.Range("AF3").Value = 0.07199999999999999
$ws.Range("AG3").Value = 0.07199999999999999
$ws.Range("AH3").Value = 0.01119402985074627
$ws.Range("AI3").Value = 0.02635431918008784
$ws.Range("AJ3").Value = 0.01119402985074627
$ws.Range("AK3").Value = 0.02635431918008784
$ws.Range("AL3").Value = 0.008
$ws.Range("AM3").Value = 0.008
$ws.Range("AN3").Value = 1.241379310344827
$ws.Range("AO3").Value = 7.25
$ws.Range("AP3").Value = 1.241379310344827
$ws.Range("AQ3").Value = 7.25

# --- Remove row 4 entirely ---
$ws.Rows(4).Delete()
